# Adds the "createPO" worksheet (Create Purchase Order test data) as the
# last sheet in the workbook, mirroring the layout/style conventions of the
# existing "campaign" sheet.

$wb = $excel.ActiveWorkbook
$campaign = $wb.Worksheets.Item(1)

# --- Create the new worksheet, positioned after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "createPO"

# --- Header row (row 1) ---
$ws.Range("A1").Value = "testID"
$ws.Range("B1").Value = "testNAme"
$ws.Range("C1").Value = "subject"
$ws.Range("D1").Value = "vendor"
$ws.Range("E1").Value = "Billing Address"
$ws.Range("F1").Value = "Shipping Address"
$ws.Range("G1").Value = "Item"
$ws.Range("H1").Value = "qty"
$ws.Range("I1").Value = "price"
$ws.Range("J1").Value = "expected msg"

# --- Data row (row 2) ---
$ws.Range("A2").Value = "tc_01"
$ws.Range("B2").Value = "createPoTest"
$ws.Range("C2").Value = "laptop"
$ws.Range("D2").Value = "ABB"
$ws.Range("E2").Value = "Bengaluru"
$ws.Range("F2").Value = "Bengaluru"
$ws.Range("G2").Value = "micro phone"
$ws.Range("H2").Value = "'3"
$ws.Range("I2").Value = "'100"
$ws.Range("J2").Value = " Purchase Order Information "

# --- Reuse the existing header/body styles from the "campaign" sheet so the
#     style table doesn't get needlessly duplicated ---
$campaign.Range("A1").Copy()
$ws.Range("A1:H1").PasteSpecial(-4122)

$campaign.Range("A2").Copy()
$ws.Range("A2:H5").PasteSpecial(-4122)

# --- I1:J1 header cells: same yellow fill as the rest of the header, but
#     only left/right borders ---
$ws.Range("I1:J1").Interior.Color = 65535
$ws.Range("I1:J1").Borders.Item(7).LineStyle = 1
$ws.Range("I1:J1").Borders.Item(10).LineStyle = 1

$excel.CutCopyMode = $false

# --- Sheet view / selection ---
$ws.Activate()
$ws.Range("J2").Select()

# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 27.22
$ws.Columns.Item(3).ColumnWidth = 13.44
$ws.Columns.Item(5).ColumnWidth = 12.89
$ws.Columns.Item(6).ColumnWidth = 15.78
$ws.Columns.Item(7).ColumnWidth = 11.55
$ws.Columns.Item(10).ColumnWidth = 25.55
